# Enable answering polls functionality
# - Populate the "answer" counts (column D) on the Answer sheet
# - Update sheet selections / active-cell bookmarks left by the author
# - Make "Answer" the active/selected sheet (it was "User" before)
# - Add a portrait page-setup to the "User" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Answer sheet: fill in the previously-zeroed "count" column (D)
# ---------------------------------------------------------------------
$wsAnswer = $wb.Worksheets.Item("Answer")

$wsAnswer.Range("D3").Value  = 3
$wsAnswer.Range("D5").Value  = 15
$wsAnswer.Range("D6").Value  = 6
$wsAnswer.Range("D7").Value  = 9
$wsAnswer.Range("D8").Value  = 12
$wsAnswer.Range("D9").Value  = 16
$wsAnswer.Range("D10").Value = 17
$wsAnswer.Range("D11").Value = 19
$wsAnswer.Range("D12").Value = 22
$wsAnswer.Range("D13").Value = 10
$wsAnswer.Range("D14").Value = 24
$wsAnswer.Range("D15").Value = 27
$wsAnswer.Range("D16").Value = 4
$wsAnswer.Range("D17").Value = 13
$wsAnswer.Range("D18").Value = 20
$wsAnswer.Range("D19").Value = 28
$wsAnswer.Range("D21").Value = 23

# ---------------------------------------------------------------------
# Update the remembered selection / active-cell on the sheets the
# author clicked through while testing the feature.
# ---------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Activate() | Out-Null
$wsUser.Range("D2").Select() | Out-Null
$wsUser.PageSetup.Orientation = 1

$wsPoll = $wb.Worksheets.Item("Poll")
$wsPoll.Activate() | Out-Null
$wsPoll.Range("C4").Select() | Out-Null

$wsQuestion = $wb.Worksheets.Item("Question")
$wsQuestion.Activate() | Out-Null
$wsQuestion.Range("B5").Select() | Out-Null

$wsOption = $wb.Worksheets.Item("Option")
$wsOption.Activate() | Out-Null
$wsOption.Range("B22").Select() | Out-Null

# Answer ends up the active sheet/selection when the author saved.
$wsAnswer.Activate() | Out-Null
$wsAnswer.Range("D21").Select() | Out-Null
